$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 4, shifting all existing
# data rows (old 4..33) down to (new 6..35). Excel's Insert copies the
# formatting (e.g. the date-format style on column D) from the row
# being pushed down, matching the style seen on the rest of the table.
$ws.Range("A4:T5").EntireRow.Insert()

# New row 4: Higo, Primera, week of 2022-05-04
$ws.Range("A4").Value = 6
$ws.Range("B4").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C4").Value = "Metropolitana"
$ws.Range("D4").Value = 44685
$ws.Range("E4").Value = 13
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100101
$ws.Range("H4").Value = "Berries"
$ws.Range("I4").Value = 100101006
$ws.Range("J4").Value = "Higo"
$ws.Range("K4").Value = "Sin especificar"
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 15000
$ws.Range("O4").Value = 15000
$ws.Range("P4").Value = 15000
$ws.Range("Q4").Value = "`$/bandeja 7 kilos"
$ws.Range("R4").Value = "Región Metropolitana"
$ws.Range("S4").Value = 2143
$ws.Range("T4").Value = 7

# New row 5: Higo, Segunda, week of 2022-05-04
$ws.Range("A5").Value = 6
$ws.Range("B5").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C5").Value = "Metropolitana"
$ws.Range("D5").Value = 44685
$ws.Range("E5").Value = 13
$ws.Range("F5").Value = "Fruta"
$ws.Range("G5").Value = 100101
$ws.Range("H5").Value = "Berries"
$ws.Range("I5").Value = 100101006
$ws.Range("J5").Value = "Higo"
$ws.Range("K5").Value = "Sin especificar"
$ws.Range("L5").Value = "Segunda"
$ws.Range("M5").Value = 70
$ws.Range("N5").Value = 12000
$ws.Range("O5").Value = 12000
$ws.Range("P5").Value = 12000
$ws.Range("Q5").Value = "`$/bandeja 7 kilos"
$ws.Range("R5").Value = "Región Metropolitana"
$ws.Range("S5").Value = 1714
$ws.Range("T5").Value = 7
